$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1199.5
$ws.Range("I43").Value = 900
$ws.Range("J43").Value = 1499
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 1499
$ws.Range("M43").Value = -831
$ws.Range("N43").Value = -1637
$ws.Range("H98").Value = 1211.8937
$ws.Range("I98").Value = 1093.186
$ws.Range("K98").Value = 1093.186
$ws.Range("M98").Value = 404.8140000000001
$ws.Range("H100").Value = 2972.6667
$ws.Range("I100").Value = 1801.2858
$ws.Range("K100").Value = 1801.2858
$ws.Range("M100").Value = -1260.2858
$ws.Range("H122").Value = 1211.8937
$ws.Range("I122").Value = 1093.186
$ws.Range("K122").Value = 3279.558
$ws.Range("M122").Value = -829.558
$ws.Range("H137").Value = 3709.842
$ws.Range("I137").Value = 3198.739
$ws.Range("K137").Value = 9596.217000000001
$ws.Range("M137").Value = -7046.217000000001
$ws.Range("H138").Value = 2966.532
$ws.Range("I138").Value = 2260.8948
$ws.Range("J138").Value = 3445.3572
$ws.Range("K138").Value = 6782.6844
$ws.Range("L138").Value = 10336.0716
$ws.Range("M138").Value = -1642.6844
$ws.Range("N138").Value = -20616.0716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1457.5428
$ws.Range("I2").Value = 902.9231
$ws.Range("K2").Value = 902.9231
$ws.Range("M2").Value = -789.9231
$ws.Range("H32").Value = 12629642
$ws.Range("J32").Value = 27785778
$ws.Range("L32").Value = 27785778
$ws.Range("N32").Value = -27786352
$ws.Range("H61").Value = 1865.1
$ws.Range("I61").Value = 1771.8975
$ws.Range("K61").Value = 1771.8975
$ws.Range("M61").Value = -1559.8975
$ws.Range("H74").Value = 3924.75
$ws.Range("J74").Value = 6699
$ws.Range("L74").Value = 6699
$ws.Range("N74").Value = -8447
$ws.Range("H77").Value = 3924.75
$ws.Range("J77").Value = 6699
$ws.Range("L77").Value = 33495
$ws.Range("N77").Value = -42231
$ws.Range("H88").Value = 4594.5
$ws.Range("I88").Value = 1168.6666
$ws.Range("J88").Value = 6650
$ws.Range("K88").Value = 1168.6666
$ws.Range("L88").Value = 6650
$ws.Range("M88").Value = -762.6666
$ws.Range("N88").Value = -7462
$ws.Range("H91").Value = 4594.5
$ws.Range("I91").Value = 1168.6666
$ws.Range("J91").Value = 6650
$ws.Range("K91").Value = 1168.6666
$ws.Range("L91").Value = 6650
$ws.Range("M91").Value = 235.3334
$ws.Range("N91").Value = -9458
$ws.Range("H97").Value = 1112.2693
$ws.Range("I97").Value = 770.4545000000001
$ws.Range("K97").Value = 770.4545000000001
$ws.Range("M97").Value = -274.4545000000001
$ws.Range("H116").Value = 1457.5428
$ws.Range("I116").Value = 902.9231
$ws.Range("K116").Value = 902.9231
$ws.Range("M116").Value = 1391.0769
$ws.Range("H128").Value = 59995
$ws.Range("J128").Value = 59995
$ws.Range("L128").Value = 59995
$ws.Range("N128").Value = -69955
$ws.Range("H136").Value = 1865.1
$ws.Range("I136").Value = 1771.8975
$ws.Range("K136").Value = 5315.6925
$ws.Range("M136").Value = -2765.6925

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1457.5428
$ws.Range("I3").Value = 902.9231
$ws.Range("K3").Value = 902.9231
$ws.Range("M3").Value = -788.9231
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 10000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 10000
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -9706
$ws.Range("H60").Value = 94050.8
$ws.Range("J60").Value = 94050.8
$ws.Range("L60").Value = 94050.8
$ws.Range("N60").Value = -95248.8
$ws.Range("H107").Value = 1880.0625
$ws.Range("I107").Value = 1759.4615
$ws.Range("K107").Value = 1759.4615
$ws.Range("M107").Value = 160.5385000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3991.4517
$ws.Range("I31").Value = 2257.2666
$ws.Range("K31").Value = 2257.2666
$ws.Range("M31").Value = -1962.2666
$ws.Range("H34").Value = 3991.4517
$ws.Range("I34").Value = 2257.2666
$ws.Range("K34").Value = 2257.2666
$ws.Range("M34").Value = -2055.2666
$ws.Range("H58").Value = 2226.1562
$ws.Range("I58").Value = 1883.4642
$ws.Range("K58").Value = 1883.4642
$ws.Range("M58").Value = -1680.4642
$ws.Range("H75").Value = 104642.664
$ws.Range("J75").Value = 104642.664
$ws.Range("L75").Value = 104642.664
$ws.Range("N75").Value = -106638.664
$ws.Range("H78").Value = 104642.664
$ws.Range("J78").Value = 104642.664
$ws.Range("L78").Value = 313927.992
$ws.Range("N78").Value = -323911.992
$ws.Range("H99").Value = 1644.5714
$ws.Range("I99").Value = 1588.6666
$ws.Range("K99").Value = 1588.6666
$ws.Range("M99").Value = -90.66660000000002
$ws.Range("H100").Value = 110780
$ws.Range("J100").Value = 110780
$ws.Range("L100").Value = 110780
$ws.Range("N100").Value = -112944
$ws.Range("H126").Value = 1644.5714
$ws.Range("I126").Value = 1588.6666
$ws.Range("K126").Value = 4765.9998
$ws.Range("M126").Value = -2295.9998
$ws.Range("H136").Value = 2226.1562
$ws.Range("I136").Value = 1883.4642
$ws.Range("K136").Value = 5650.392599999999
$ws.Range("M136").Value = -3100.392599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 44.75
$ws.Range("I6").Value = 46.333332
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 138.999996
$ws.Range("L6").Value = 120
$ws.Range("M6").Value = -25.99999600000001
$ws.Range("N6").Value = -346
$ws.Range("H12").Value = 71.5
$ws.Range("I12").Value = 27.5
$ws.Range("K12").Value = 82.5
$ws.Range("M12").Value = 90.5
$ws.Range("H107").Value = 574.4706
$ws.Range("J107").Value = 594.2
$ws.Range("L107").Value = 1782.6
$ws.Range("N107").Value = -5622.6
$ws.Range("H131").Value = 1705
$ws.Range("I131").Value = 1239
$ws.Range("J131").Value = 1832.091
$ws.Range("K131").Value = 3717
$ws.Range("L131").Value = 5496.272999999999
$ws.Range("M131").Value = 1323
$ws.Range("N131").Value = -15576.273

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 134749
$ws.Range("J128").Value = 134749
$ws.Range("L128").Value = 134749
$ws.Range("N128").Value = -144709
$ws.Range("H132").Value = 3841.0557
$ws.Range("I132").Value = 3946.5
$ws.Range("J132").Value = 2997.5
$ws.Range("K132").Value = 11839.5
$ws.Range("L132").Value = 8992.5
$ws.Range("M132").Value = -9309.5
$ws.Range("N132").Value = -14052.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3268.3684
$ws.Range("I82").Value = 2686.375
$ws.Range("K82").Value = 2686.375
$ws.Range("M82").Value = -2325.375
$ws.Range("H85").Value = 3268.3684
$ws.Range("I85").Value = 2686.375
$ws.Range("K85").Value = 2686.375
$ws.Range("M85").Value = -1438.375
$ws.Range("H110").Value = 75481.75
$ws.Range("J110").Value = 75481.75
$ws.Range("L110").Value = 75481.75
$ws.Range("N110").Value = -83661.75
$ws.Range("H122").Value = 2940.5417
$ws.Range("I122").Value = 2917.5217
$ws.Range("J122").Value = 3470
$ws.Range("K122").Value = 8752.5651
$ws.Range("L122").Value = 10410
$ws.Range("M122").Value = -6302.5651
$ws.Range("N122").Value = -15310
$ws.Range("H136").Value = 3428.6453
$ws.Range("I136").Value = 2534.9546
$ws.Range("K136").Value = 7604.8638
$ws.Range("M136").Value = -5054.8638

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 58338.5
$ws.Range("J46").Value = 58338.5
$ws.Range("L46").Value = 58338.5
$ws.Range("N46").Value = -58800.5
$ws.Range("H64").Value = 99981.75
$ws.Range("J64").Value = 99981.75
$ws.Range("L64").Value = 99981.75
$ws.Range("N64").Value = -100477.75
$ws.Range("H67").Value = 99981.75
$ws.Range("J67").Value = 99981.75
$ws.Range("L67").Value = 99981.75
$ws.Range("N67").Value = -101697.75
$ws.Range("H134").Value = 58338.5
$ws.Range("J134").Value = 58338.5
$ws.Range("L134").Value = 175015.5
$ws.Range("N134").Value = -180085.5
$ws.Range("H136").Value = 17425.258
$ws.Range("I136").Value = 1386.8636
$ws.Range("J136").Value = 49502.047
$ws.Range("K136").Value = 4160.5908
$ws.Range("L136").Value = 148506.141
$ws.Range("M136").Value = -1610.5908
$ws.Range("N136").Value = -153606.141
